$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing rows 2-5 down to 3-6.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A2").Value = "parth"
$ws.Range("B2").Value = "parthpatel082828@gmail.com"
$ws.Range("C2").Value = 45406.67672707176
$ws.Range("D2").Value = "15:52:4"
$ws.Range("E2").Value = "16:12:36"

# Re-apply the same date number format already used by the other rows in
# column C (numFmtId 14), reusing the existing style instead of creating one.
$ws.Range("C3").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C2").Value = 45406.67672707176
